$d = $word.ActiveDocument

# 1. Append a new paragraph at the very end of the document (after the two
#    trailing empty paragraphs) containing the new "Anmerkung" remark.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "Anmerkung 23.09.2020: Registrierung sollte einfacher erfolgen -> NFC-Sticker. Terminverwaltung nicht nötig. "

# 2. Move the "_GoBack" bookmark from its old location (end of the
#    "...bertragungen nachverfolgen zu können. " paragraph) onto the very
#    end of the freshly added paragraph (collapsed, right after the text,
#    before the paragraph mark).
#
#    A collapsed Range positioned exactly at a paragraph-mark boundary
#    confuses Bookmarks.Add in this runtime, so we work around it: insert a
#    throwaway placeholder character right after the target spot, anchor the
#    bookmark just before it (now a safe, non-boundary position), then
#    remove the placeholder again. The bookmark stays collapsed in place.
$endPos = $newPara.Range.End - 1
$placeholderRange = $d.Range($endPos, $endPos)
$placeholderRange.InsertAfter("X")

$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholder = $d.Range($endPos, $endPos + 1)
$placeholder.Delete()
